$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to be treated as text so that numeric-looking strings
# (e.g. "1.001", "4.479") are preserved exactly instead of being converted
# to floating point numbers by Excels automatic type inference.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '26.491.91'
$ws.Range("E2").Value = '  +1.77%  '
$ws.Range("D3").Value = '1.671.41'
$ws.Range("E3").Value = '  +1.49%  '
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = '219.75'
$ws.Range("E5").Value = '  +2.19%  '
$ws.Range("E6").Value = '  +1.17%  '
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("E8").Value = '  +2.69%  '
$ws.Range("D9").Value = '0.06376'
$ws.Range("E9").Value = '  +0.11%  '
$ws.Range("D10").Value = '21.78'
$ws.Range("E10").Value = '  +4.59%  '
$ws.Range("D11").Value = '0.07799'
$ws.Range("E11").Value = '  +1.52%  '
$ws.Range("D12").Value = '1.679.26'
$ws.Range("E12").Value = '  +1.95%  '
$ws.Range("D13").Value = '4.479'
$ws.Range("E13").Value = '  +1.18%  '
$ws.Range("D14").Value = '0.5576'
$ws.Range("E14").Value = '  +0.58%  '
$ws.Range("D15").Value = '0.0₅8307'
$ws.Range("E15").Value = '  -0.32%  '
$ws.Range("D16").Value = '65.56'
$ws.Range("D17").Value = '26.493.50'
$ws.Range("E17").Value = '  +1.74%  '
$ws.Range("E18").Value = '  +0.02%  '
$ws.Range("D19").Value = '4.763'
$ws.Range("E19").Value = '  +1.00%  '
$ws.Range("E20").Value = '  +2.54%  '
$ws.Range("D21").Value = '10.32'
$ws.Range("E21").Value = '  +1.56%  '
$ws.Range("D22").Value = '6.305'
$ws.Range("E22").Value = '  +0.67%  '
$ws.Range("D23").Value = '1.002'
$ws.Range("E23").Value = '  +0.09%  '
$ws.Range("D24").Value = '0.1270'
$ws.Range("D25").Value = '138.54'
$ws.Range("E25").Value = '  -4.33%  '
$ws.Range("D26").Value = '7.390'
$ws.Range("E26").Value = '  -0.19%  '
$ws.Range("D27").Value = '16.30'
$ws.Range("E27").Value = '  +2.50%  '
$ws.Range("D28").Value = '1.426'
$ws.Range("E28").Value = '  +2.78%  '
$ws.Range("D29").Value = '0.06217'
$ws.Range("E29").Value = '  +4.17%  '
$ws.Range("E30").Value = '  +1.91%  '
$ws.Range("D31").Value = '3.609'
$ws.Range("E31").Value = '  +5.94%  '
$ws.Range("E32").Value = '  +0.55%  '
$ws.Range("E33").Value = '  +2.12%  '
$ws.Range("E34").Value = '  +1.09%  '
$ws.Range("D35").Value = '0.6107'
$ws.Range("E35").Value = '  +8.46%  '
$ws.Range("D36").Value = '2.415'
$ws.Range("E36").Value = '  +0.92%  '
$ws.Range("D37").Value = '2.781'
$ws.Range("E37").Value = '  +1.00%  '
$ws.Range("D39").Value = '6.037'
$ws.Range("E39").Value = '  +3.04%  '
$ws.Range("D40").Value = '1.091.38'
$ws.Range("E40").Value = '  +6.25%  '
$ws.Range("D41").Value = '0.8564'
$ws.Range("E41").Value = '  +0.23%  '
$ws.Range("E42").Value = '  -0.01%  '
$ws.Range("D43").Value = '100.55'
$ws.Range("E43").Value = '  +1.76%  '
$ws.Range("D44").Value = '1.816.74'
$ws.Range("E44").Value = '  +1.17%  '
$ws.Range("B45").Value = 'BabyDogeCoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D45").Value = '0.0₈112'
$ws.Range("E45").Value = '  +1.66%  '
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").Value = '58.27'
$ws.Range("E46").Value = '  +4.61%  '
$ws.Range("E47").Value = '  -0.50%  '
$ws.Range("D48").Value = '1.519'
$ws.Range("E48").Value = '  +10.25%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").Value = '8.090'
$ws.Range("E49").Value = '  +0.45%  '
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").Value = '0.05193'
$ws.Range("E50").Value = '  +0.95%  '
$ws.Range("B51").Value = 'Aptos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D51").Value = '6.003'
$ws.Range("E51").Value = '  +1.35%  '
